# fix: remove isGlobalAdmin field in user (#2235)
#
# The "user" worksheet has one column per user field. The
# "is_global_admin" field (column AH, between "is_admin" and
# "is_forbidden") is removed entirely, shifting all later columns
# (is_forbidden, is_deleted, signup_application, created_ip,
# last_signin_time, last_signin_ip) one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the "is_global_admin" header column dynamically so the script
# is resilient to minor layout differences; fall back to the known
# column letter (AH) if the header can't be located for some reason.
$targetCol = $null
$headerRow = $ws.Rows.Item(1)
$found = $headerRow.Find("is_global_admin")
if ($found -ne $null) {
    $targetCol = $found.Column
}

if ($targetCol -ne $null) {
    $ws.Columns.Item($targetCol).EntireColumn.Delete()
} else {
    $ws.Columns("AH:AH").Delete()
}

# Leave the active selection where it would naturally land after
# deleting the column (on the cell now occupied by the shifted
# "signup_application" value in the first data row).
$ws.Range("AK2").Select()
